$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the defined names "rihanna3" -> "jamilawoods1" (both scoped
#    names, one per sheet: Sheet1!rihanna3 and Sheet3!rihanna3).
# ---------------------------------------------------------------------
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    $n.Name = "jamilawoods1"
}

# ---------------------------------------------------------------------
# 2. New tracklist data (Jamila Woods - HEAVN) replacing the old
#    Rihanna "ANTI" tracklist. Columns: B=Title, C=Composer, D=Performer,
#    E=Time (fraction of a day).
# ---------------------------------------------------------------------

function Set-Tracklist($ws) {
    $ws.Range("B2").Value = "Bubbles"
    $ws.Range("C2").Value = "Jamila Woods"
    $ws.Range("D2").Value = "Jamila Woods"
    $ws.Range("E2").Value = 0.087500000000000008

    $ws.Range("B3").Value = "VRY BLK feat. Noname"
    $ws.Range("C3").Value = "Jamila Woods, Noname"
    $ws.Range("D3").Value = "Jamila Woods, Noname"
    $ws.Range("E3").Value = 0.14791666666666667

    $ws.Range("B4").Value = "Lonely Lonely feat. Lorine Chia"
    $ws.Range("C4").Value = "Jamila Woods"
    $ws.Range("D4").Value = "Jamila Woods, Lorine Chia"
    $ws.Range("E4").Value = 0.14305555555555557

    $ws.Range("B5").Value = "HEAVN"
    $ws.Range("C5").ClearContents()
    $ws.Range("D5").Value = "Jamila Woods"
    $ws.Range("E5").Value = 0.15277777777777776

    $ws.Range("B6").Value = "In My Name"
    $ws.Range("C6").Value = "Jamila Woods"
    $ws.Range("D6").Value = "Jamila Woods"
    $ws.Range("E6").Value = 0.075694444444444439

    $ws.Range("B7").Value = "LSD feat. Chance The Rapper"
    $ws.Range("C7").Value = "Chance The Rapper, Jamila Woods"
    $ws.Range("D7").Value = "Chance The Rapper, Jamila Woods"
    $ws.Range("E7").Value = 0.16388888888888889

    $ws.Range("B8").Value = "Emerald St. feat. Saba"
    $ws.Range("C8").Value = "Jamila Woods"
    $ws.Range("D8").Value = "Jamila Woods"
    $ws.Range("E8").Value = 0.14097222222222222

    $ws.Range("B9").Value = "Walk My Way"
    $ws.Range("C9").Value = "Jamila Woods, Saba"
    $ws.Range("D9").Value = "Jamila Woods, Saba"
    $ws.Range("E9").Value = 0.13402777777777777

    $ws.Range("B10").Value = "Lately"
    $ws.Range("C10").Value = "Jamila Woods"
    $ws.Range("D10").Value = "Jamila Woods"
    $ws.Range("E10").Value = 0.13055555555555556

    $ws.Range("B11").Value = "Breadcrumbs feat. Donnie Trumpet"
    $ws.Range("C11").Value = "Jamila Woods"
    $ws.Range("D11").Value = "Jamila Woods, Donnie Trumpet"
    $ws.Range("E11").Value = 0.16874999999999998

    $ws.Range("B12").Value = "Stellar"
    $ws.Range("C12").Value = "Jamila Woods"
    $ws.Range("D12").Value = "Jamila Woods"
    $ws.Range("E12").Value = 0.084722222222222213

    $ws.Range("B13").Value = "Holy"
    $ws.Range("C13").Value = "Jamila Woods"
    $ws.Range("D13").Value = "Jamila Woods"
    $ws.Range("E13").Value = 0.13333333333333333

    $ws.Range("B14").Value = "Way Up"
    $ws.Range("C14").Value = "Jamila Woods"
    $ws.Range("D14").Value = "Jamila Woods"
    $ws.Range("E14").Value = 0.16388888888888889

    # Match the slightly wider columns used to fit the new (shorter but
    # differently distributed) text.
    $ws.Columns.Item(2).ColumnWidth = 32.6640625
    $ws.Columns.Item(3).ColumnWidth = 31.88671875
    $ws.Columns.Item(4).ColumnWidth = 31.88671875

    # E18 picks up the same number formatting as its neighbours once the
    # sheet is refreshed/re-saved.
    $ws.Range("E18").NumberFormat = $ws.Range("F18").NumberFormat
}

Set-Tracklist($wb.Worksheets.Item("Sheet1"))
Set-Tracklist($wb.Worksheets.Item("Sheet3"))

# ---------------------------------------------------------------------
# 3. Force a full recalculation so that the formula-driven "preview"
#    sheet (Sheet2), which mirrors Sheet1 via LEFTB/REPT formulas, picks
#    up the new values and cached results.
# ---------------------------------------------------------------------
$excel.CalculateFullRebuild()
